$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string / comment text for the new journal entry
$newText = "J'ai améliorer la vérification des coordonnées"

# Grow the table (Tableau1) by one row so the new data becomes part of it
$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.ListRows.Add() | Out-Null

# Add the new row 43 data (table Tableau1 currently spans B4:L42)
$ws.Range("B43").Value = 44267
$ws.Range("C43").Value = 0.5625
$ws.Range("D43").Value = 0.57291666666666663
$ws.Range("E43").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure fin]]),"""",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure début]])"
$ws.Range("F43").Value = "Ma-20"
$ws.Range("G43").Value = "Code"
$ws.Range("H43").Value = "Jeu"
$ws.Range("I43").Value = "CPNV"
$ws.Range("J43").Value = $newText
$ws.Range("K43").Value = "Oui"
$ws.Range("L43").Value = "LVT"

# Update selection / view position to match the saved workbook state
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("H47").Select()
